# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-24 08:22:04
#
# Updates the Y2 GIT & Liver session-analysis sheet: a few "Recorded By"
# lists get re-ordered, the MICROBIOLOGY C1 session #1 (row 12) has now
# been recorded (so its row flips from "Not Recorded" to "Recorded" and
# picks up the recorder + attendance figures), and the dependent
# statistics (recorded/missing session counts, coverage %, average
# attendance %) are refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ANATOMY C1 #1): reorder "Recorded By" list -----------------
$ws.Cells.Item(2, 7).Value = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, System"

# --- Row 3 (ANATOMY C1 #2): reorder "Recorded By" list -----------------
$ws.Cells.Item(3, 7).Value = "hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# --- Row 4 (ANATOMY C1 #3): reorder "Recorded By" list -----------------
$ws.Cells.Item(4, 7).Value = "hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# --- Row 7 (BIOCHEMISTRY LAB/CBL C1 #1): reorder "Recorded By" list ----
$ws.Cells.Item(7, 7).Value = "menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg"

# --- Row 12 (MICROBIOLOGY C1 #1): session has now been recorded --------
$ws.Cells.Item(12, 7).Value = "yassmina.fattoh@med.asu.edu.eg"
$ws.Cells.Item(12, 8).Value = "24/251"
$ws.Cells.Item(12, 9).Value = "Recorded"

# Re-colour row 12 from the "Not Recorded" (pink) look to the "Recorded"
# (green) look by copying the format of an already-recorded row.
$recordedFormat = $ws.Range("A2:I2")
$recordedFormat.Copy()
$ws.Range("A12:I12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Class Statistics panel (K/L column, ANATOMY block) ----------------
# Recorded Sessions: 11 -> 12
$ws.Cells.Item(6, 12).Value = 12
# Missing Sessions: 2 -> 1
$ws.Cells.Item(7, 12).Value = 1
# Coverage %: 37.9% -> 41.4%  (leading ' keeps these as literal text, not
# an auto-converted numeric percentage, matching the original cells)
$ws.Cells.Item(9, 12).Value = "'41.4%"
# Average Attendance %: 24.9% -> 23.6%
$ws.Cells.Item(10, 12).Value = "'23.6%"

# --- Group Statistics panel (row 15: Year 2 / C1) -----------------------
# Recorded: 11 -> 12
$ws.Cells.Item(15, 15).Value = 12
# Missing: 2 -> 1
$ws.Cells.Item(15, 16).Value = 1
# Coverage %: 37.9% -> 41.4%
$ws.Cells.Item(15, 18).Value = "'41.4%"
# Avg Attendance %: 24.9% -> 23.6%
$ws.Cells.Item(15, 19).Value = "'23.6%"

# The leading apostrophes above mark the cells "text", which also stamps a
# quote-prefix flag onto their style. Restore each cell's original look by
# re-pasting formats from a same-styled neighbour so only the text changed.
$statFmtSrc = $ws.Range("K9")         # s=4, same style as L9/L10
$statFmtSrc.Copy()
$ws.Range("L9").PasteSpecial(-4122)
$statFmtSrc.Copy()
$ws.Range("L10").PasteSpecial(-4122)

$groupFmtSrc = $ws.Range("N15")       # s=4, same style as O15/P15/R15/S15
$groupFmtSrc.Copy()
$ws.Range("R15").PasteSpecial(-4122)
$groupFmtSrc.Copy()
$ws.Range("S15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
